{"js": "// The fixture's \"CaseText\" paragraphs were authored with embedded/trailing\n// literal newlines inside their <w:t> runs (an artifact of how the fixture\n// text was wrapped when it was generated). This cleans that up:\n//   - a trailing newline at the end of a paragraph's text is dropped\n//   - newlines in the middle of a paragraph's text are collapsed to a\n//     single space, joining the wrapped lines back into one line\n// Paragraphs that don't contain a literal \"\\n\" (e.g. normal single-line\n// paragraphs, TOC fields, etc.) are left completely untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load('items/style');\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load('text');\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n\n  // Only touch the body-copy paragraphs (style \"Case Text\") that actually\n  // contain embedded newline characters in their run text.\n  if (paragraph.style !== 'Case Text') continue;\n  const text = paragraph.text;\n  if (text.indexOf('\\n') === -1) continue;\n\n  const fixed = text.replace(/\\s*\\n\\s*/g, ' ').replace(/\\s+$/, '');\n  if (fixed !== text) {\n    paragraph.insertText(fixed, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The fixture's \"Case Text\" paragraphs were authored with embedded/trailing\n# literal newlines inside their run text (an artifact of how the fixture\n# text was wrapped when it was generated). This cleans that up:\n#   - a trailing newline at the end of a paragraph's text is dropped\n#   - newlines in the middle of a paragraph's text are collapsed to a\n#     single space, joining the wrapped lines back into one line\n# Paragraphs that don't contain a literal newline (e.g. normal single-line\n# paragraphs, TOC fields, etc.) are left completely untouched.\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $p = $paras.Item($i)\n    if ($p.Style.NameLocal -ne \"Case Text\") { continue }\n\n    $full = $p.Range.Text\n    # Range.Text includes the trailing paragraph mark (CR); strip it off\n    # before inspecting/rewriting the run text, then it's implicitly\n    # restored because we only replace the Range's text content.\n    $body = $full.TrimEnd([char]13)\n\n    if ($body.IndexOf([char]10) -lt 0) { continue }\n\n    $fixed = [System.Text.RegularExpressions.Regex]::Replace($body, \"\\s*\\n\\s*\", \" \")\n    $fixed = $fixed.TrimEnd()\n\n    if ($fixed -ne $body) {\n        $p.Range.Text = $fixed\n    }\n}\n\nWrite-Output \"done\"\n"}
